# Hien thong bao them the loai thanh cong
#
# The edit:
#   1. The trailing paragraph that previously held the hidden "_GoBack"
#      bookmark (a lone space run, right after the controller paragraphs)
#      loses that bookmark.
#   2. Two new one-cell tables (plus a couple of connecting paragraphs) are
#      appended just before the document's final empty paragraph. The first
#      table holds the CateController "getAdd"/"postAdd"/"getList" PHP code
#      (the "_GoBack" bookmark now lives inside it, around the
#      `return view('admin.cate.list');` line). The second table holds the
#      Blade snippet that renders the flash/success message.

$d = $word.ActiveDocument

# --- 1. Drop the hidden _GoBack bookmark from its old home -----------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 2. Build the replacement OOXML fragment (2 tables + paragraphs) -------
$newContentXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="9350"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9350" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>public function getAdd() {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/><w:t>return view(''admin.cate.add'');</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    }</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">    public function postAdd(CateRequest $request) {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>$cate = new Cate;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">$cate-&gt;name </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= $request-&gt;txtCateName;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">$cate-&gt;alias </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= $request-&gt;txtCateName;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">$cate-&gt;order </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= $request-&gt;txtOrder;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">$cate-&gt;parent_id </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= 1;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">$cate-&gt;keywords </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= $request-&gt;txtKeywords;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">$cate-&gt;description </w:t></w:r><w:r><w:tab/><w:t>= $request-&gt;txtDescription;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>$cate-&gt;save();</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>return redirect()-&gt;route(''admin.cate.list'')-&gt;</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>with([''level_message''=&gt;''success'' ,''flash_message''=&gt;''Success'']);</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    }</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">    public function getList() {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:tab/><w:t>return view(''admin.cate.list'');</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t xml:space="preserve">    }</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Hiển thị thông báo thêm thành công</w:t></w:r></w:p><w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="9350"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9350" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>&lt;div class="col-lg-12"&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">                        @if (Session::has(''flash_message''))</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">                            &lt;div class="alert alert-{!! Session::get(''level_message'') !!}"&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">                                {!! Session::get(''flash_message'') !!}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">                            &lt;/div&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">                        @endif</w:t></w:r></w:p><w:p><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">                    &lt;/div&gt;</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'

# --- 3. Insert it right before the document's final (empty) paragraph ------
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($paraCount)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.InsertXML($newContentXml)
